$d = $word.ActiveDocument

# 1) Merge "DAY 0" + "4" runs into a single "DAY 04" paragraph.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "DAY 04"

# 2) The paragraph that follows (previously empty) gets the date text,
#    centered, bold, size 14pt (w:sz 28 half-points), matching the
#    paragraph-mark run formatting already present on that paragraph.
$p2 = $d.Paragraphs.Item(2)
$p2.Alignment = 1          # wdAlignParagraphCenter
$r2 = $p2.Range
$r2.Text = "22.03.2023"
$r2.Font.Bold = 1
$r2.Font.Size = 14

Write-Output "done"
